$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of data (row 74) after the existing last row (row 73).
$row = 74

# Column A: store the date as plain text (matches existing inline string rows),
# not as a numeric date serial.
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "10/30/2025"

# Column B: numeric profit value.
$ws.Cells.Item($row, 2).Value = 10998.9
